# Apply crypto price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "42.676.55"
$ws.Cells.Item(2, 5).Value = "  -0.87%  "
$ws.Cells.Item(3, 4).Value = "2.510.87"
$ws.Cells.Item(3, 5).Value = "  -1.78%  "
$ws.Cells.Item(4, 5).Value = "  +0.10%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "317.37"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +4.09%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "95.77"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -2.95%  "
$ws.Cells.Item(7, 5).Value = "  +0.54%  "
$ws.Cells.Item(8, 5).Value = "  +0.00%  "
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.537"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -2.21%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "36.08"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -2.55%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.0809"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -0.72%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "7.58"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -2.43%  "
$ws.Cells.Item(13, 5).Value = "  -2.57%  "
$ws.Cells.Item(14, 4).Value = "2.901.64"
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "15.51"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +4.27%  "
$ws.Cells.Item(16, 4).Value = "2.511.78"
$ws.Cells.Item(16, 5).Value = "  -0.97%  "
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "0.857"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -3.15%  "
$ws.Cells.Item(18, 4).Value = "42.681.55"
$ws.Cells.Item(18, 5).Value = "  -1.01%  "
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "13.10"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -4.57%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0970"
$ws.Cells.Item(20, 5).Value = "  -2.00%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "6.56"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -1.13%  "
$ws.Cells.Item(22, 5).Value = "  -0.91%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "251.82"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -2.07%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "2.99"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.09%  "
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "2.04"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -3.11%  "
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "26.97"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -3.88%  "
$ws.Cells.Item(27, 5).Value = "  +0.13%  "
$ws.Cells.Item(28, 5).Value = "  +11.89%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "38.79"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +1.73%  "
$ws.Cells.Item(30, 5).Value = "  -0.91%  "
$ws.Cells.Item(31, 5).Value = "  -2.52%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "155.93"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -1.94%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "19.43"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +2.42%  "
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "3.35"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +0.50%  "
$ws.Cells.Item(35, 5).Value = "  -3.52%  "
$ws.Cells.Item(36, 5).Value = "  -2.57%  "
$ws.Cells.Item(37, 5).Value = "  -5.39%  "
$ws.Cells.Item(38, 5).Value = "  -2.23%  "
$ws.Cells.Item(39, 2).Value = "Stellar"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.120"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +0.20%  "
$ws.Cells.Item(40, 2).Value = "EnergySwap"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "24.01"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -7.43%  "
$ws.Cells.Item(41, 5).Value = "  +3.04%  "
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "3.85"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.95%  "
$ws.Cells.Item(43, 5).Value = "  -2.24%  "
$ws.Cells.Item(44, 5).Value = "  +0.13%  "
$ws.Cells.Item(45, 5).Value = "  -1.98%  "
$ws.Cells.Item(46, 4).Value = "2.043.93"
$ws.Cells.Item(46, 5).Value = "  -2.34%  "
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "84.51"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -2.93%  "
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "8.84"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -2.27%  "
$ws.Cells.Item(49, 4).Value = "2.755.11"
$ws.Cells.Item(49, 5).Value = "  -1.75%  "
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "73.39"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -3.05%  "
$ws.Cells.Item(51, 5).Value = "  -1.09%  "
